# Add daily power records
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Fill in existing rows 120 and 121 (Date/Start Time/End Time)
$ws.Range("A120").Value = 43443
$ws.Range("B120").Value = 0
$ws.Range("C120").Value = 0

$ws.Range("A121").Value = 43444
$ws.Range("B121").Value = 0
$ws.Range("C121").Value = 0

# Extend the table by 4 rows (122-125)
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null

# Row 122
$ws.Range("A122").Value = 43445
$ws.Range("B122").Value = 0.62152777777777779
$ws.Range("C122").Value = 0.65
$ws.Range("D122").Formula = "=(C122-B122)* 1440"
$ws.Range("E122").Formula = "=IF(C122>B122, (C122-B122)*1440, (B122-C122)*1440)"
$ws.Range("F122").Formula = "=ABS((C122-B122)*1440)"

# Row 123
$ws.Range("A123").Value = 43445
$ws.Range("B123").Value = 0.65763888888888888
$ws.Range("C123").Value = 0.73263888888888884
$ws.Range("D123").Formula = "=(C123-B123)* 1440"
$ws.Range("E123").Formula = "=IF(C123>B123, (C123-B123)*1440, (B123-C123)*1440)"
$ws.Range("F123").Formula = "=ABS((C123-B123)*1440)"

# Row 124 (no Date/Start/End, just calculated columns)
$ws.Range("D124").Formula = "=(C124-B124)* 1440"
$ws.Range("E124").Formula = "=IF(C124>B124, (C124-B124)*1440, (B124-C124)*1440)"
$ws.Range("F124").Formula = "=ABS((C124-B124)*1440)"

# Row 125 (no Date/Start/End, just calculated columns)
$ws.Range("D125").Formula = "=(C125-B125)* 1440"
$ws.Range("E125").Formula = "=IF(C125>B125, (C125-B125)*1440, (B125-C125)*1440)"
$ws.Range("F125").Formula = "=ABS((C125-B125)*1440)"

# Update view/selection to mirror the authored edit
$win = $excel.ActiveWindow
$win.ScrollRow = 109
$win.ScrollColumn = 1
$ws.Range("B125").Select() | Out-Null

Write-Host "Done. Table range:" $lo.Range.Address()
